$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the "VALOR MORA" total and the "Cant. Periodos" count
# ------------------------------------------------------------------
$ws.Range("E11").Value = 448841
$ws.Range("F13").Value = 3

# ------------------------------------------------------------------
# 2. Swap the "Novedad de Ingreso" / "Novedad de Retiro" column headers
# ------------------------------------------------------------------
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"

# ------------------------------------------------------------------
# 3. Insert two new data rows right after the existing detail rows
#    (old rows 19-24 -- blank rows + signature block -- shift down to 21-26)
# ------------------------------------------------------------------
$ws.Rows("19:20").Insert()

# Preserve the old "closing" (bottom-border) row format + values of the
# last detail row (row 18, WILLIAN ESTEBAN / period 2508) onto the new
# final row (20) before we restyle row 18.
$ws.Range("B18:J18").Copy($ws.Range("B20:J20"))

# Row 18 is no longer the last row of the table, so give it the regular
# "middle" row styling (same as row 17), keeping its own data values.
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143328807"
$ws.Range("D18").Value = "WILLIAN ESTEBAN PAJARO BANQUEZ"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 82708
$ws.Range("G18").Value = 2067698

# New row 19: JHON JAIRO POLO CASTRO, period 2509 (middle-row styling)
$ws.Range("B17:J17").Copy($ws.Range("B19:J19"))
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143327340"
$ws.Range("D19").Value = "JHON JAIRO POLO CASTRO"
$ws.Range("E19").Value = "2509"
$ws.Range("F19").Value = 94475
$ws.Range("G19").Value = 2361859

# New row 20: WILLIAN ESTEBAN PAJARO BANQUEZ, period 2509 (closing-row styling,
# values/format already copied from the old row 18 above)
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143328807"
$ws.Range("D20").Value = "WILLIAN ESTEBAN PAJARO BANQUEZ"
$ws.Range("E20").Value = "2509"
$ws.Range("F20").Value = 82708
$ws.Range("G20").Value = 2067698

# The "Periodo Mora" column is now centered for every detail row
$ws.Range("E16:E20").HorizontalAlignment = -4108
